$wb = $excel.ActiveWorkbook

# Sheet: VENTAS POR GRUPO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M7").Value = 629.3200000000001

# Sheet: VENTA MENSUAL
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F7").Value = 2228.69
$ws2.Range("F22").Value = 10261.6

# Sheet: CUMPLIMIENTO MENSUAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D16").Value = 7307.46
$ws3.Range("E16").Value = 36958.78
$ws3.Range("F16").Value = 0.1650797537807593
$ws3.Range("D19").Value = 10261.6
$ws3.Range("E19").Value = 55116.39762291768
$ws3.Range("F19").Value = 0.1569580038101823
